# Commit: [TEST SCRAPE] updated files from azure vm
#
# 1) "ODI Batting" sheet: clear the stray empty INNING_NUMBER (col B) cells
#    on the rows where the player did not bat (rows 5, 8, 9, 10).
# 2) Add a new sheet "ODI Batting Extra" (after "ODI Bowling") with per-match
#    batting extras (batting position, 4s, 6s, % of team runs, man of match).

$wb = $excel.ActiveWorkbook

# --- 1) Clear stray empty cells in "ODI Batting" ---------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B5").ClearContents()
$batting.Range("B8").ClearContents()
$batting.Range("B9").ClearContents()
$batting.Range("B10").ClearContents()

# --- 2) Add the new "ODI Batting Extra" sheet -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Columns A, C, D, E hold text-like data (match codes, counts, percentages
# formatted as strings) -- force Text format so they are not reinterpreted
# as numbers. Column B (batting position) and F (MAN_OF_MATCH "NO") keep
# their natural types (number / text respectively).
$extra.Range("A1:A17").NumberFormat = "@"
$extra.Range("C1:E17").NumberFormat = "@"

# Header row
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $extra.Cells.Item(1, $col)
    $cell.Value2 = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$rows = @(
    @("4605", 6,    "3", "1", "12.46%", "NO"),
    @("4608", 6,    "2", "1", "7.31%",  "NO"),
    @("4614", 6,    "5", "2", "13.06%", "NO"),
    @("4625", 6,    $null, $null, $null, "NO"),
    @("4639", 7,    "1", "0", "4.25%",  "NO"),
    @("4649", 6,    "2", "2", "19.42%", "NO"),
    @("4669", 6,    $null, $null, $null, "NO"),
    @("4673", $null,$null, $null, $null, "NO"),
    @("4676", 6,    $null, $null, $null, "NO"),
    @("4686", $null,$null, $null, $null, "NO"),
    @("4688", 6,    "0", "0", "1.15%",  "NO"),
    @("4690", $null,$null, $null, $null, "NO"),
    @("4692", $null,$null, $null, $null, "NO"),
    @("4695", 6,    "5", "0", "33.33%", "NO"),
    @("4697", 6,    "0", "0", "1.69%",  "NO"),
    @("4735", 6,    "3", "2", "14.23%", "NO")
)

$r = 2
foreach ($row in $rows) {
    $extra.Cells.Item($r, 1).Value2 = $row[0]
    if ($null -ne $row[1]) {
        $extra.Cells.Item($r, 2).Value2 = $row[1]
    }
    if ($null -ne $row[2]) {
        $extra.Cells.Item($r, 3).Value2 = $row[2]
    }
    if ($null -ne $row[3]) {
        $extra.Cells.Item($r, 4).Value2 = $row[3]
    }
    if ($null -ne $row[4]) {
        $extra.Cells.Item($r, 5).Value2 = $row[4]
    }
    $extra.Cells.Item($r, 6).Value2 = $row[5]
    $r++
}

Write-Output "edit complete"
